$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 4 data rows (A1:AA4). A new (mostly blank) record
# is being appended as row 5, and the used range grows one column further
# right (to AB) to hold it. Insert the new row/column first so nothing
# existing shifts unexpectedly.
$ws.Rows(5).Insert()
$ws.Columns(28).Insert()

# Row 4 had a stray, value-less "Memo" cell (U4) left over from data entry;
# the edit drops it entirely.
$ws.Range("U4").ClearContents()

# Populate the new row 5: Ward = "100A" and the checkbox-style column (R) is
# unchecked (FALSE); the rest of the row stays blank except for a touch in
# column AB so the sheet's used range reports through AB5.
$ws.Range("A5").Value = "100A"
$ws.Range("R5").Value = $false
$ws.Range("AB5").Value = " "
